$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds numbers stored as text; force Text format before
# assigning so Excel doesn't silently coerce the strings back to numbers.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1008"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "129"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "114"

$ws.Range("B5").Value = "Kendrick Lamar"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "105"

$ws.Range("B6").Value = "Kids See Ghosts"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "98"

$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "94"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "73"

$ws.Range("B10").Value = "JAY-Z"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "73"

$ws.Range("B11").Value = "Frank Ocean"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "72"
